$d = $word.ActiveDocument

# Pair 0: paragraph 1
$p = $d.Paragraphs.Item(1)
$r = $p.Range
$res = $r.Find.Execute("Introspektivni pogled na odnos između Hawkinga i kontingenta prostora / vremena. Ovaj film iznosi galilejske i newtonske zakone i postoji odnos s Einsteinovom teorijom opće relativnosti. Film je metodički režiran, izlaže detalje o čovjeku (Hawking) kao i njegovu djelu (Crne rupe). Intervjui s njegovom obitelji malo su predugi pa je nažalost ", $true, $false, $false, $false, $false, $true, 1, $false, "Introspektivni pogled na odnos između Hawkinga i kontingenta prostora / vremena. ## Ovaj film iznosi galilejske i newtonske zakone i postoji odnos s Einsteinovom teorijom opće relativnosti. ## Film je metodički režiran, izlaže detalje o čovjeku (Hawking) kao i njegovu djelu (Crne rupe). ## Intervjui s njegovom obitelji malo su predugi pa je nažalost ", 2)
Write-Host "Pair 0 (para 1):" $res

# Pair 1: paragraph 1
$p = $d.Paragraphs.Item(1)
$r = $p.Range
$res = $r.Find.Execute(". Soundtrack Philipa Glassa izvrsno nadopunjuje film. Samo je jedan drugi čovjek mogao skladati takve jezive zvjezdane melodije (Jean Michel Jarre). Sve u svemu, toplo bih preporučio ovaj film na osnovu Hawkingovih", $true, $false, $false, $false, $false, $true, 1, $false, ". ## Soundtrack Philipa Glassa izvrsno nadopunjuje film. ## Samo je jedan drugi čovjek mogao skladati takve jezive zvjezdane melodije (Jean Michel Jarre). ## Sve u svemu, toplo bih preporučio ovaj film na osnovu Hawkingovih", 2)
Write-Host "Pair 1 (para 1):" $res

# Pair 2: paragraph 3
$p = $d.Paragraphs.Item(3)
$r = $p.Range
$res = $r.Find.Execute(" tako grozan film dugo ... dugo vremena ... ", $true, $false, $false, $false, $false, $true, 1, $false, " tako grozan film dugo ... dugo vremena ...  ## ", 2)
Write-Host "Pair 2 (para 3):" $res

# Pair 3: paragraph 3
$p = $d.Paragraphs.Item(3)
$r = $p.Range
$res = $r.Find.Execute(" otići nakon 20 minuta ... ", $true, $false, $false, $false, $false, $true, 1, $false, " otići nakon 20 minuta ... ## ", 2)
Write-Host "Pair 3 (para 3):" $res

# Pair 4: paragraph 3
$p = $d.Paragraphs.Item(3)
$r = $p.Range
$res = $r.Find.Execute(" ulogu ... ", $true, $false, $false, $false, $false, $true, 1, $false, " ulogu ... ## ", 2)
Write-Host "Pair 4 (para 3):" $res

# Pair 5: paragraph 3
$p = $d.Paragraphs.Item(3)
$r = $p.Range
$res = $r.Find.Execute(" glume? sudeći po ", $true, $false, $false, $false, $false, $true, 1, $false, " glume? ## sudeći po ", 2)
Write-Host "Pair 5 (para 3):" $res

# Pair 6: paragraph 3
$p = $d.Paragraphs.Item(3)
$r = $p.Range
$res = $r.Find.Execute(" u bliskoj budućnosti ... oboje izgledaju prilično lijepo .. ", $true, $false, $false, $false, $false, $true, 1, $false, " u bliskoj budućnosti ... ## oboje izgledaju prilično lijepo .. ## ", 2)
Write-Host "Pair 6 (para 3):" $res

# Pair 7: paragraph 6
$p = $d.Paragraphs.Item(6)
$r = $p.Range
$res = $r.Find.Execute("na ovaj proizvod kad su mi rekli da sam anemična. Sad ih uzimam otprilike 4 mjeseca i anemije više nema. Dobar proizvod. Lako se probavlja (za razliku od nekih drugih dodataka željezu).", $true, $false, $false, $false, $false, $true, 1, $false, "na ovaj proizvod kad su mi rekli da sam anemična. ## Sad ih uzimam otprilike 4 mjeseca i anemije više nema. ## Dobar proizvod. ## Lako se probavlja (za razliku od nekih drugih dodataka željezu).", 2)
Write-Host "Pair 7 (para 6):" $res

# Pair 8: paragraph 9
$p = $d.Paragraphs.Item(9)
$r = $p.Range
$res = $r.Find.Execute("Ovo je jedan od mojih najdražih deserta i brzo se topi u ustima. Ova marka je dobra i isporučuje ", $true, $false, $false, $false, $false, $true, 1, $false, "Ovo je jedan od mojih najdražih deserta i brzo se topi u ustima. ## Ova marka je dobra i isporučuje ", 2)
Write-Host "Pair 8 (para 9):" $res

# Pair 9: paragraph 9
$p = $d.Paragraphs.Item(9)
$r = $p.Range
$res = $r.Find.Execute(". Svi bi trebali jednom ", $true, $false, $false, $false, $false, $true, 1, $false, ". ## Svi bi trebali jednom ", 2)
Write-Host "Pair 9 (para 9):" $res

# Pair 10: paragraph 9
$p = $d.Paragraphs.Item(9)
$r = $p.Range
$res = $r.Find.Execute(". ", $true, $false, $false, $false, $false, $true, 1, $false, ". ## ", 2)
Write-Host "Pair 10 (para 9):" $res

# Pair 11: paragraph 12
$p = $d.Paragraphs.Item(12)
$r = $p.Range
$res = $r.Find.Execute("Ovo je fantastična zagonetka / poklon za mlade I stare. Riječ je o 32 trokutasta jaka magnetska komada koji se mogu sastaviti na više različitih načina. Jednostavno je super i imat ćete problema držati ga podalje od odraslih.", $true, $false, $false, $false, $false, $true, 1, $false, "Ovo je fantastična zagonetka / poklon za mlade I stare. ## Riječ je o 32 trokutasta jaka magnetska komada koji se mogu sastaviti na više različitih načina. ## Jednostavno je super i imat ćete problema držati ga podalje od odraslih.", 2)
Write-Host "Pair 11 (para 12):" $res

# Pair 12: paragraph 15
$p = $d.Paragraphs.Item(15)
$r = $p.Range
$res = $r.Find.Execute(" je još jedan loš zombi film. U usporedbi s većinom ostalih, jedina razlika ovdje je ženska osoba. Zaplet je isti. Akcijske scene nisu privlačne. Posebni efekti su ", $true, $false, $false, $false, $false, $true, 1, $false, " je još jedan loš zombi film. ## U usporedbi s većinom ostalih, jedina razlika ovdje je ženska osoba. ## Zaplet je isti. ## Akcijske scene nisu privlačne. ## Posebni efekti su ", 2)
Write-Host "Pair 12 (para 15):" $res

# Pair 13: paragraph 18
$p = $d.Paragraphs.Item(18)
$r = $p.Range
$res = $r.Find.Execute("YoYo izgleda izvan ravnoteže. Bez obzira na ", $true, $false, $false, $false, $false, $true, 1, $false, "YoYo izgleda izvan ravnoteže. ## Bez obzira na ", 2)
Write-Host "Pair 13 (para 18):" $res

# Pair 14: paragraph 18
$p = $d.Paragraphs.Item(18)
$r = $p.Range
$res = $r.Find.Execute(" je naginjati se na jednu stranu. ", $true, $false, $false, $false, $false, $true, 1, $false, " je naginjati se na jednu stranu. ## ", 2)
Write-Host "Pair 14 (para 18):" $res

# Pair 15: paragraph 18
$p = $d.Paragraphs.Item(18)
$r = $p.Range
$res = $r.Find.Execute(". Imam malo iskustva s ", $true, $false, $false, $false, $false, $true, 1, $false, ". ## Imam malo iskustva s ", 2)
Write-Host "Pair 15 (para 18):" $res

# Pair 16: paragraph 21
$p = $d.Paragraphs.Item(21)
$r = $p.Range
$res = $r.Find.Execute(". ", $true, $false, $false, $false, $false, $true, 1, $false, ". ## ", 2)
Write-Host "Pair 16 (para 21):" $res

# Pair 17: paragraph 21
$p = $d.Paragraphs.Item(21)
$r = $p.Range
$res = $r.Find.Execute(". Oni jednostavno neće zaštititi / podržavati vaše ", $true, $false, $false, $false, $false, $true, 1, $false, ". ## Oni jednostavno neće zaštititi / podržavati vaše ", 2)
Write-Host "Pair 17 (para 21):" $res
